# Apply edits to the "Data Path C" worksheet (sheet2):
#  1. Rename the SH column header from SH[4:0] to SH[3:0]
#  2. Shrink the placeholder values in that column from 5 Ø's to 4 Ø's
#  3. Change the frozen panes: freeze at C3 (columns B-and-left + rows 1-2)
#     instead of just freezing at column C, and update the pane selections.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Path C")

# --- 1 & 2: update column N (SH[...]) header + data ---
$ws.Range("N2").Value = "SH[3:0]"
$ws.Range("N3").Value = " ØØØØ"
$ws.Range("N4").Value = "ØØØØ"
$ws.Range("N5").Value = "ØØØØ"
$ws.Range("N6").Value = "ØØØØ"
$ws.Range("N7").Value = "ØØØØ"
$ws.Range("N8").Value = "ØØØØ"
$ws.Range("N9").Value = "ØØØØ"
$ws.Range("N10").Value = "ØØØØ"
$ws.Range("N11").Value = "ØØØØ"
$ws.Range("N12").Value = "ØØØØ"
$ws.Range("N13").Value = "0011"

# --- 3: re-freeze panes at C3 ---
$ws.Activate() | Out-Null
$ws.Range("C3").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("C1").Select() | Out-Null
$ws.Range("A3").Select() | Out-Null
$ws.Range("L13").Select() | Out-Null
